$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Date value
$ws.Range("B8").Value = "2024-07-15T11:25:06-04:00"

# Update Context values (strip the IG-specific URL prefix, keep element:Resource)
$ws.Range("B21").Value = "element:CarePlan"
$ws.Range("B22").Value = "element:ServiceRequest"
